$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$t = $ws.ListObjects.Item(2)
Write-Host "Table Name: $($t.Name) Range before: $($t.Range.Address())"

# delete the raw column
$ws.Columns.Item(13).Delete()

Write-Host "Table Range after col delete: $($t.Range.Address())"

try {
    $t.Resize($ws.Range("B4:N43"))
    Write-Host "Resize succeeded"
} catch {
    Write-Host "Resize failed: $_"
}
Write-Host "Table Range after resize: $($t.Range.Address())"
